$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Wins (AD1), Losses (AE1), Ties (AF1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold / bordered / centered header style used by the rest of row 1 (e.g. AC1)
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Season record for every player row: Wins=71, Losses=91, Ties=0
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 71
    $ws.Cells.Item($r, 31).Value = 91
    $ws.Cells.Item($r, 32).Value = 0
}
